# Add new function "tottus" data: update existing rows 2-13 with refreshed
# order info and append 3 new data rows (14-16) to the OrdenSalida sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrdenSalida")

function Set-TextCell($sheet, $row, $col, $val) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

function Set-NumberCell($sheet, $row, $col, $val) {
    $sheet.Cells.Item($row, $col).Value = $val
}

# Column indices:
# A=1 B=2 C=3 D=4 E=5 F=6 G=7 H=8 S=19 T=20 AB=28

$rowsData = @{
    2  = @{ A = "47035710"; B = "47035710"; G = "20250814"; H = "20250826"; S = "20215631"; T = 24;  AB = "8865" }
    3  = @{ A = "47035710"; B = "47035710"; G = "20250814"; H = "20250826"; S = "20287252"; T = 24;  AB = "8865" }
    4  = @{ A = "47035710"; B = "47035710"; G = "20250814"; H = "20250826"; S = "20287256"; T = 48;  AB = "8865" }
    5  = @{ A = "47035710"; B = "47035710"; G = "20250814"; H = "20250826"; S = "20287253"; T = 24;  AB = "8865" }
    6  = @{ A = "47035710"; B = "47035710"; G = "20250814"; H = "20250826"; S = "20202318"; T = 24;  AB = "8865" }
    7  = @{ A = "47035710"; B = "47035710"; G = "20250814"; H = "20250826"; S = "20215632"; T = 240; AB = "8865" }
    8  = @{ A = "47035710"; B = "47035710"; G = "20250814"; H = "20250826"; S = "20202336"; T = 24;  AB = "8865" }
    9  = @{ A = "47035710"; B = "47035710"; G = "20250814"; H = "20250826"; S = "20202335"; T = 12;  AB = "8865" }
    10 = @{ A = "47035710"; B = "47035710"; G = "20250814"; H = "20250826"; S = "20202339"; T = 72;  AB = "8865" }
    11 = @{ A = "47035710"; B = "47035710"; G = "20250814"; H = "20250826"; S = "20202309"; T = 60;  AB = "8865" }
    12 = @{ A = "47035710"; B = "47035710"; G = "20250814"; H = "20250826"; S = "20202351"; T = 24;  AB = "8865" }
    13 = @{ A = "47035710"; B = "47035710"; G = "20250814"; H = "20250826"; S = "20202337"; T = 12;  AB = "8865" }
    14 = @{ A = "47035710"; B = "47035710"; C = "78627210-6"; D = "Hipermercados TOTTUS SA"; E = "CD TOTTUS"; F = "CD TOTTUS"; G = "20250814"; H = "20250826"; S = "20287251"; T = 24; AB = "8865" }
    15 = @{ A = "47035710"; B = "47035710"; C = "78627210-6"; D = "Hipermercados TOTTUS SA"; E = "CD TOTTUS"; F = "CD TOTTUS"; G = "20250814"; H = "20250826"; S = "20202310"; T = 12; AB = "8865" }
    16 = @{ A = "47035710"; B = "47035710"; C = "78627210-6"; D = "Hipermercados TOTTUS SA"; E = "CD TOTTUS"; F = "CD TOTTUS"; G = "20250814"; H = "20250826"; S = "20215634"; T = 48; AB = "8865" }
}

$textCols = @{ A = 1; B = 2; C = 3; D = 4; E = 5; F = 6; G = 7; H = 8; S = 19; AB = 28 }
$numCols  = @{ T = 20 }

foreach ($rowNum in 2..16) {
    $rowVals = $rowsData[$rowNum]
    foreach ($key in $rowVals.Keys) {
        $value = $rowVals[$key]
        if ($textCols.ContainsKey($key)) {
            Set-TextCell $ws $rowNum $textCols[$key] $value
        } elseif ($numCols.ContainsKey($key)) {
            Set-NumberCell $ws $rowNum $numCols[$key] $value
        }
    }
}

Write-Output "Updated rows 2-16 on sheet $($ws.Name)"
